$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.194.30'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '2.643.85'
$ws.Range("E3").Value = '  -1.94%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''527.59'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '''144.76'
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.569'
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("D9").Value = '''6.66'
$ws.Range("E9").Value = '  -4.88%  '
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '3.108.02'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").Value = '59.125.46'
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").Value = '''21.02'
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = '2.666.56'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '''342.09'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = '''6.35'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '''65.57'
$ws.Range("E23").Value = '  +2.93%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = '''0.168'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").Value = '''0.996'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '''7.26'
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = '0.0₃0800'
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("E29").Value = '  -4.79%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = '''18.97'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '''150.13'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("D34").Value = '''4.23'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''1.21'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").Value = '''0.929'
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D37").Value = '''0.877'
$ws.Range("E37").Value = '  -2.65%  '
$ws.Range("D38").Value = '''1.49'
$ws.Range("E38").Value = '  -2.33%  '
$ws.Range("D39").Value = '''36.60'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").Value = '''0.997'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  -5.99%  '
$ws.Range("D43").Value = '''0.0974'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").Value = '''272.26'
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("D45").Value = '''19.45'
$ws.Range("E45").Value = '  -4.00%  '
$ws.Range("D46").Value = '''0.0540'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").Value = '2.053.24'
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("D49").Value = '''4.82'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''19.16'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '''0.0230'
$ws.Range("E51").Value = '  -1.41%  '
